# Recalibrated extrapolation model results after removing option quotes
# priced under USD 5 (treated as noise) from the calibration inputs.
# Only the fitted output columns (ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN)
# for the affected expiries change; all other data is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - 2025-11-28
$ws.Range("D4").Value = 115588.2525109894
$ws.Range("E4").Value = -0.01224077990359004
$ws.Range("F4").Value = 0.179229836183176
$ws.Range("G4").Value = -1.004376605563255
$ws.Range("H4").Value = 9.568966754044327

# Row 5 - 2025-12-26
$ws.Range("D5").Value = 116210.3841605729
$ws.Range("E5").Value = -0.01876016353704567
$ws.Range("F5").Value = 0.2110870876821966
$ws.Range("G5").Value = -1.412266168921447
$ws.Range("H5").Value = 11.92861027360767

# Row 6 - 2026-01-30
$ws.Range("D6").Value = 117077.6657309065
$ws.Range("E6").Value = -0.01873107783446396
$ws.Range("F6").Value = 0.1989536075674057
$ws.Range("G6").Value = -0.6000217679317379
$ws.Range("H6").Value = 5.907066210377918

# Row 7 - 2026-02-27
$ws.Range("D7").Value = 117643.4294263742
$ws.Range("E7").Value = -0.02221224002388503
$ws.Range("F7").Value = 0.2035558849628692
$ws.Range("G7").Value = -0.5735490706616849
$ws.Range("H7").Value = 5.486543673453048

# Row 8 - 2026-05-29
$ws.Range("D8").Value = 117943.7718681236
$ws.Range("E8").Value = -0.03342676312100754
$ws.Range("F8").Value = 0.1899021976984366
$ws.Range("G8").Value = -1.224330677113425
$ws.Range("H8").Value = 8.765372173076505

# Row 10 - 2026-09-25
$ws.Range("D10").Value = 121104.4003866674
$ws.Range("E10").Value = -0.09787097527198921
$ws.Range("F10").Value = 0.3925788658842621
$ws.Range("G10").Value = -1.782893227396202
$ws.Range("H10").Value = 9.397937649256962

# Row 13 - 2025-09-16
$ws.Range("D13").Value = 114028.2038096776
$ws.Range("E13").Value = -0.01297389370140632
$ws.Range("F13").Value = 0.1160858698519495
$ws.Range("G13").Value = -0.2925775517246894
$ws.Range("H13").Value = 4.9482245951439

# Row 14 - 2025-09-17
$ws.Range("D14").Value = 113972.4301232762
$ws.Range("E14").Value = -0.01808311698640697
$ws.Range("F14").Value = 0.1364557410892673
$ws.Range("G14").Value = -0.7546068212074712
$ws.Range("H14").Value = 8.359587161463592

# Row 17 - 2025-09-12
$ws.Range("D17").Value = 113892.031510227
$ws.Range("E17").Value = -0.03117937887729277
$ws.Range("F17").Value = 0.1485199320274635
$ws.Range("G17").Value = -0.5434261775285498
$ws.Range("H17").Value = 3.89298725166621

# Row 19 - 2025-10-03
$ws.Range("D19").Value = 114936.5091921413
$ws.Range("E19").Value = 0.009237655493657
$ws.Range("F19").Value = 0.1515487264717648
$ws.Range("G19").Value = -0.7894635930751819
$ws.Range("H19").Value = 6.664206672859858

# Row 20 - 2025-10-10
$ws.Range("D20").Value = 114929.8521208075
$ws.Range("E20").Value = 0.003273882903961198
$ws.Range("F20").Value = 0.1500736356555293
$ws.Range("G20").Value = -0.6080044081389298
$ws.Range("H20").Value = 5.557220421447647
